# "Challenge initialisation and test in main"
# Adds a second challenge row ("Defi ton voisin de droite ...") to the
# Feuil1 sheet, mirroring the formatting already used for the first data
# row, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New challenge row (row 3) -------------------------------------------
$ws.Range("A3").Value = "Defi ton voisin de droite au pierre feuille ciseau, celui qui pert bois 3 gorgées"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "[soirée chill, soirée bar]"
$ws.Range("D3").Value = 1

# --- Formatting -------------------------------------------------------------
# The whole data columns use wrapped text; re-asserting it here also folds
# the header cell D1 (which previously had its own "no-wrap" variant of the
# header style) onto the same wrapped style as A1:C1, and stamps the new
# row and the previously-unstyled D2 cell with the wrap style too.
$ws.Columns.Item(1).WrapText = $true
$ws.Columns.Item(2).WrapText = $true
$ws.Columns.Item(3).WrapText = $true
$ws.Columns.Item(4).WrapText = $true

# New row is tall enough to show the full challenge text wrapped.
$ws.Rows.Item(3).RowHeight = 51

# --- View state ---------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 720
$win.Width = 29400
$win.Height = 18400

$ws.Range("F4").Select()
